# Adds two new columns, I ("I0") and J ("IF"), with per-row data, to the
# single worksheet of the workbook. Mirrors the existing H ("IP") column's
# header style and extends the sheet's used-range dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold font, thin border, centered/top alignment)
# from the existing "IP" header cell (H1) onto the two new header cells,
# matching the s="1" styling used by the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2-48) ----------------------------------------------
$data = @(
    @(2,8,9),
    @(3,8,8),
    @(4,7,8),
    @(5,7,8),
    @(6,8,8),
    @(7,9,9),
    @(8,7,7),
    @(9,7,7),
    @(10,8,8),
    @(11,8,9),
    @(12,9,9),
    @(13,8,8),
    @(14,9,10),
    @(15,8,8),
    @(16,8,9),
    @(17,9,9),
    @(18,8,8),
    @(19,9,9),
    @(20,10,10),
    @(21,8,9),
    @(22,8,8),
    @(23,9,9),
    @(24,9,9),
    @(25,6,7),
    @(26,7,7),
    @(27,9,9),
    @(28,7,7),
    @(29,8,8),
    @(30,7,7),
    @(31,5,6),
    @(32,6,7),
    @(33,6,7),
    @(34,11,11),
    @(35,5,5),
    @(36,10,10),
    @(37,5,6),
    @(38,9,9),
    @(39,8,8),
    @(40,5,6),
    @(41,4,4),
    @(42,8,8),
    @(43,7,7),
    @(44,5,6),
    @(45,6,6),
    @(46,9,9),
    @(47,5,5),
    @(48,6,6)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
